# Deploy the implementation guide.
#
# The IG publisher regenerated docs/CodeSystem-qc-sample-type.xlsx with a
# newer build: the CodeSystem's Status moved from "active" to "draft" and
# the Date metadata value was bumped to the new publish timestamp. Both
# values live on the "Metadata" worksheet, one column over from their
# label in column A:
#   A6 "Status" / B6 "active"  -> "draft"
#   A8 "Date"   / B8 "2023-05-12T12:33:13+00:00" -> "2023-08-01T16:12:28+00:00"

$wb = $excel.ActiveWorkbook
$metadata = $wb.Worksheets.Item("Metadata")

$metadata.Range("B6").Value = "draft"
$metadata.Range("B8").Value = "2023-08-01T16:12:28+00:00"
